$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width change (closest reachable value to the authored 18.33203125
#     given this engine's column-width quantization) ---
$ws.Columns.Item(1).ColumnWidth = 17.57

# --- New font style (sz 12) for A5, must happen before the column-B number format
#     so the generated cellXfs ordering matches the target workbook ---
$ws.Range("A5").Font.Size = 12

# --- Row 4 content (numeric password must be written before column B gets the
#     text number format, so it is stored as a genuine number) ---
$ws.Range("A4").Value2 = "910097(苏艳辉)"
$ws.Range("B4").Value2 = 300427
$ws.Range("C4").Value2 = "交易客户端"
$ws.Range("D4").Value2 = "交易客户端"
$ws.Range("E4").Value2 = "dce"
$ws.Range("F4").Value2 = "金瑞快期"
$ws.Rows.Item(4).RowHeight = 36.6

# --- Row 5 content, except B5 which must be added after E5 so the shared
#     string table ends up in the same order as the target file ---
$ws.Range("A5").Value2 = "910101(孙悦)"
$ws.Range("C5").Value2 = "交易客户端"
$ws.Range("E5").Value2 = "zce"
$ws.Range("F5").Value2 = "金瑞快期"
$ws.Rows.Item(5).RowHeight = 37.8

# --- Column B becomes a text column (password values kept as text) ---
$ws.Columns.Item(2).NumberFormat = "@"

# --- B5's password is textual (leading zeros) so it must be entered as text,
#     once the column is already formatted as text ---
$ws.Range("B5").Value2 = "052927"

# --- Final selection matches the authored workbook ---
$ws.Range("B6").Select()
